$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are plain-text price strings (e.g. "53.878.70", "0.999")
# that must stay text. Force text format so Excel does not auto-convert
# numeric-looking strings to numbers (which would drop trailing zeros),
# then restore the default "Normal" style so no formatting diff is left
# behind on the cell.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '54.045.83'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.257.29'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('E4').Value = '  -0.82%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '495.12'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '128.51'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.83%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.524'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.99%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0949'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('E10').Value = '  +0.88%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.336'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.74'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.35%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.653.61'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.68'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '53.999.04'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.71%  '
$ws.Range('E16').Value = '  +0.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.247.83'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.22'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.14'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '300.46'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.40%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.29'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.60%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '60.70'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.99%  '
$ws.Range('E24').Value = '  +0.17%  '
$ws.Range('E25').Value = '  -1.37%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.27'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '170.46'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.84%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.60'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0688'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.91'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.47%  '
$ws.Range('E31').Value = '  +0.89%  '
$ws.Range('E32').Value = '  -0.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '17.75'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.997'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.70%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.943'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.19'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.70'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.88%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.372'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.39'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.35'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '125.07'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.36%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.78'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.37%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0491'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0890'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.542'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.27%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '239.44'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.43%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.371'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.59%  '
$ws.Range('E48').Value = '  +0.25%  '
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.09'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.50%  '
$ws.Range('E51').Value = '  -0.84%  '
